# Splits three long single-run paragraphs ("Programa" PT/EN text and the
# "Bibliografia" text) into multiple <w:t> runs separated by manual line
# breaks (<w:br/>), matching how the document was re-wrapped upstream.
#
# Strategy: use Find & Replace with the "^l" special character (manual
# line break) in the replacement text. Word inserts a real line-break
# element and splits the text run around it, which produces exactly the
# <w:t>...</w:t><w:br/><w:t>...</w:t> pattern required by the target
# OOXML.

$d = $word.ActiveDocument

function Insert-LineBreakAfter($findText) {
    $range = $d.Content
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, ($findText + "^l"), 2) | Out-Null
}

# Some snippets ("Eukarya.") occur once in the Portuguese paragraph and once
# in the English paragraph, so the plain text is ambiguous document-wide.
# Disambiguate by matching extra trailing context and keeping that context
# in the replacement (after the inserted break).
function Insert-LineBreakAfterWithContext($findText, $context) {
    $range = $d.Content
    $range.Find.Execute(($findText + $context), $true, $false, $false, $false, $false, $true, 1, $false, ($findText + "^l" + $context), 2) | Out-Null
}

# --- Paragraph 1: "Programa" (Portuguese) ---------------------------------
Insert-LineBreakAfter "microrganismos procarióticos e"
Insert-LineBreakAfter "Bacteria, Archaea e"
Insert-LineBreakAfterWithContext "Eukarya." "- Análise"
Insert-LineBreakAfter "microscopia ótica e eletrônica."
Insert-LineBreakAfter "via de exocitose e endocitose)."
Insert-LineBreakAfter "material genético: estrutura e função"
Insert-LineBreakAfter "mitose e meiose."

# --- Paragraph 2: "Programa" (English) ------------------------------------
Insert-LineBreakAfter "prokaryotic microorganisms and"
Insert-LineBreakAfter "Bacteria, Archaea and"
Insert-LineBreakAfterWithContext "Eukarya." "Microscope"
Insert-LineBreakAfter "optical and electron microscope."
Insert-LineBreakAfter "endocytosis and exocytosis)."
Insert-LineBreakAfter "material organization: structure and function"
Insert-LineBreakAfter "mitosis and meiosis"

# --- Paragraph 3: "Bibliografia" ------------------------------------------
Insert-LineBreakAfter "Artmed Editora Ltda, 2010."
Insert-LineBreakAfter "Artmed Editora Ltda, 3ª Edição, 2007."
Insert-LineBreakAfter "Artmed Editora, 8ª Edição, 2010."
Insert-LineBreakAfter "Editora Guanabara Koogan, 2007."
Insert-LineBreakAfter "14 Edição, 2016. "
Insert-LineBreakAfter "Editora Guanabara Koogan, 2006."
